$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row
$ws.Range("A1").Value = "Images"
$ws.Range("B1").Value = "Big_Nose"
$ws.Range("C1").Value = "Male"
$ws.Range("D1").Value = "No_Beard"

# D1 is a brand-new header cell; give it the same header formatting
# (bold font, border, centered alignment) as the existing header cells.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$data = @(
    @("data/img_align_celeba/031792.jpg", 0, 0, 1),
    @("data/img_align_celeba/126495.jpg", 0, 1, 1),
    @("data/img_align_celeba/199368.jpg", 1, 1, 0),
    @("data/img_align_celeba/046738.jpg", 0, 0, 1),
    @("data/img_align_celeba/024241.jpg", 1, 0, 1),
    @("data/img_align_celeba/098642.jpg", 0, 0, 1),
    @("data/img_align_celeba/112331.jpg", 0, 1, 0),
    @("data/img_align_celeba/071437.jpg", 0, 0, 1),
    @("data/img_align_celeba/176928.jpg", 1, 1, 0),
    @("data/img_align_celeba/061338.jpg", 0, 1, 1),
    @("data/img_align_celeba/045515.jpg", 0, 1, 0),
    @("data/img_align_celeba/038856.jpg", 0, 0, 1),
    @("data/img_align_celeba/186394.jpg", 0, 0, 1),
    @("data/img_align_celeba/147629.jpg", 0, 1, 1),
    @("data/img_align_celeba/073634.jpg", 0, 0, 1),
    @("data/img_align_celeba/135044.jpg", 1, 1, 1),
    @("data/img_align_celeba/085972.jpg", 0, 0, 1),
    @("data/img_align_celeba/171846.jpg", 1, 1, 0),
    @("data/img_align_celeba/069847.jpg", 0, 1, 1),
    @("data/img_align_celeba/198471.jpg", 0, 0, 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
